# Subnets_Info.xlsx edit: fill in calculated subnet division table (8 /27 subnets)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the frozen-pane split (target sheetView has no <pane>)
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false

# ---------------------------------------------------------------------------
# 2. Drop the extra blank rows (11-23) - target table is only 10 rows (A1:E10)
# ---------------------------------------------------------------------------
$ws.Rows("11:23").Delete()

# ---------------------------------------------------------------------------
# 3. Title text change
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Subnet IDs"

# ---------------------------------------------------------------------------
# 4. Data table values
#    Header row (row 2) text is unchanged; rows 3-10 get the calculated
#    subnet information for eight /27 subnets of 192.168.4.0/24.
# ---------------------------------------------------------------------------
$data = @(
  @("192.168.4.0",   "/27", "192.168.4.1 - 192.168.4.62",   30, "192.168.4.63"),
  @("192.168.4.64",  "/27", "192.168.4.65 - 192.168.4.126", 30, "192.168.4.127"),
  @("192.168.4.128", "/27", "192.168.4.129 - 192.168.4.190",30, "192.168.4.191"),
  @("192.168.4.192", "/27", "192.168.4.193 - 192.168.4.254",30, "192.168.4.255"),
  @("192.168.4.256", "/27", "192.168.4.257 - 192.168.4.298",30, "192.168.4.299"),
  @("192.168.4.320", "/27", "192.168.4.321 - 192.168.4.382",30, "192.168.4.383"),
  @("192.168.4.384", "/27", "192.168.4.385 - 192.168.4.446",30, "192.168.4.447"),
  @("192.168.4.448", "/27", "192.168.4.449 - 192.168.4.510",30, "192.168.4.511")
)

$r = 3
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}

Write-Output "done"
